# Update pl_mw.xlsx (Case_2_210, 380 kV) result values for rows 2-25
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.3555940612908728
$ws.Range("C2").Value = 0.09032199088997572
$ws.Range("E2").Value = 0.4204824085400674
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.2849405714051017
$ws.Range("H2").Value = 0.4728366860157251
$ws.Range("K2").Value = 0.3698186826210872
$ws.Range("O2").Value = 1.440362061093396

# Row 3
$ws.Range("B3").Value = 0.3113034052889816
$ws.Range("C3").Value = 0.08509834168479813
$ws.Range("E3").Value = 0.3669313368163643
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.2900318735753089
$ws.Range("H3").Value = 0.4791680595924674
$ws.Range("K3").Value = 0.3230553435882655
$ws.Range("O3").Value = 1.46417179861831

# Row 4
$ws.Range("B4").Value = 0.2840365355257006
$ws.Range("C4").Value = 0.08188131934259957
$ws.Range("E4").Value = 0.334125222672526
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.2934741600505184
$ws.Range("H4").Value = 0.4833271066329061
$ws.Range("K4").Value = 0.2942206238048186
$ws.Range("O4").Value = 1.480022108767201

# Row 5
$ws.Range("B5").Value = 0.2729075524543134
$ws.Range("C5").Value = 0.08056802476838243
$ws.Range("E5").Value = 0.3207738067291359
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.2949561359333188
$ws.Range("H5").Value = 0.4850901794627944
$ws.Range("K5").Value = 0.2824402993478827
$ws.Range("O5").Value = 1.486790037218604

# Row 6
$ws.Range("B6").Value = 0.271058554819831
$ws.Range("C6").Value = 0.08034981539745445
$ws.Range("E6").Value = 0.3185578214490903
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.2952069926747143
$ws.Range("H6").Value = 0.4853870554347033
$ws.Range("K6").Value = 0.2804823976044304
$ws.Range("O6").Value = 1.487932475120857

# Row 7
$ws.Range("B7").Value = 0.2838865162095203
$ws.Range("C7").Value = 0.08186361709718426
$ws.Range("E7").Value = 0.3339450924051732
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.2934938261466407
$ws.Range("H7").Value = 0.4833506078782577
$ws.Range("K7").Value = 0.2940618705370639
$ws.Range("O7").Value = 1.480112134110968

# Row 8
$ws.Range("B8").Value = 0.3403380093939745
$ws.Range("C8").Value = 0.08852294546413475
$ws.Range("E8").Value = 0.402001696597722
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.2866302458077072
$ws.Range("H8").Value = 0.4749633472463728
$ws.Range("K8").Value = 0.3537203734099137
$ws.Range("O8").Value = 1.448315704582313

# Row 9
$ws.Range("B9").Value = 0.450442963096549
$ws.Range("C9").Value = 0.1015012029118054
$ws.Range("E9").Value = 0.5361198713248569
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.275693017892813
$ws.Range("H9").Value = 0.4606731252111445
$ws.Range("K9").Value = 0.4697188147564475
$ws.Range("O9").Value = 1.395764292973837

# Row 10
$ws.Range("B10").Value = 0.5309504188337542
$ws.Range("C10").Value = 0.110982923532319
$ws.Range("E10").Value = 0.6351624008490546
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.2692132945757137
$ws.Range("H10").Value = 0.4514923991899877
$ws.Range("K10").Value = 0.5543130424975118
$ws.Range("O10").Value = 1.363175670066383

# Row 11
$ws.Range("B11").Value = 0.5674870205707521
$ws.Range("C11").Value = 0.1152839809950308
$ws.Range("E11").Value = 0.6803522854587953
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.2666072589557587
$ws.Range("H11").Value = 0.4476028064661364
$ws.Range("K11").Value = 0.5926556050264082
$ws.Range("O11").Value = 1.349667490229308

# Row 12
$ws.Range("B12").Value = 0.581309483266125
$ws.Range("C12").Value = 0.1169108330968243
$ws.Range("E12").Value = 0.6974856016378936
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.2656698697153104
$ws.Range("H12").Value = 0.4461712229172647
$ws.Range("K12").Value = 0.6071542574371165
$ws.Range("O12").Value = 1.34474244490454

# Row 13
$ws.Range("B13").Value = 0.5783331647349996
$ws.Range("C13").Value = 0.1165605460305699
$ws.Range("E13").Value = 0.6937946827214887
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.2658695481852291
$ws.Range("H13").Value = 0.4464777009476251
$ws.Range("K13").Value = 0.6040326520546841
$ws.Range("O13").Value = 1.345794667726324

# Row 14
$ws.Range("B14").Value = 0.5686244705211436
$ws.Range("C14").Value = 0.1154178610835288
$ws.Range("E14").Value = 0.6817614265349619
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.266529146210722
$ws.Range("H14").Value = 0.4474842006109938
$ws.Range("K14").Value = 0.5938488398721233
$ws.Range("O14").Value = 1.349258486602508

# Row 15
$ws.Range("B15").Value = 0.5626758811738171
$ws.Range("C15").Value = 0.1147176875169862
$ws.Range("E15").Value = 0.6743934740748045
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.2669396201241057
$ws.Range("H15").Value = 0.4481060944145057
$ws.Range("K15").Value = 0.5876082243762255
$ws.Range("O15").Value = 1.351404973000712

# Row 16
$ws.Range("B16").Value = 0.52856086287332
$ws.Range("C16").Value = 0.1107015841388517
$ws.Range("E16").Value = 0.6322119775050368
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.2693905080904386
$ws.Range("H16").Value = 0.4517523708690661
$ws.Range("K16").Value = 0.5518043894642233
$ws.Range("O16").Value = 1.364085029859282

# Row 17
$ws.Range("B17").Value = 0.5076097157763684
$ws.Range("C17").Value = 0.1082346272211652
$ws.Range("E17").Value = 0.6063706090068592
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.2709817846323048
$ws.Range("H17").Value = 0.4540627507033221
$ws.Range("K17").Value = 0.5298035651339319
$ws.Range("O17").Value = 1.372201667388666

# Row 18
$ws.Range("B18").Value = 0.4955510561069048
$ws.Range("C18").Value = 0.1068145545369958
$ws.Range("E18").Value = 0.59151994929924
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.2719291799629602
$ws.Range("H18").Value = 0.4554186067351864
$ws.Range("K18").Value = 0.5171361606404332
$ws.Range("O18").Value = 1.376994009003369

# Row 19
$ws.Range("B19").Value = 0.4914668288424195
$ws.Range("C19").Value = 0.1063335490685944
$ws.Range("E19").Value = 0.5864938921892673
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.2722554607831995
$ws.Range("H19").Value = 0.4558823087194952
$ws.Range("K19").Value = 0.5128449658643603
$ws.Range("O19").Value = 1.378637860573093

# Row 20
$ws.Range("B20").Value = 0.5098408459773793
$ws.Range("C20").Value = 0.1084973582088935
$ws.Range("E20").Value = 0.6091201508124442
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.2708090622330417
$ws.Range("H20").Value = 0.4538140131974515
$ws.Range("K20").Value = 0.5321469529189642
$ws.Range("O20").Value = 1.371324811302472

# Row 21
$ws.Range("B21").Value = 0.571476511913545
$ws.Range("C21").Value = 0.1157535468161797
$ws.Range("E21").Value = 0.6852953068051306
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.2663340612852068
$ws.Range("H21").Value = 0.4471874451530411
$ws.Range("K21").Value = 0.5968406439700971
$ws.Range("O21").Value = 1.348235909702908

# Row 22
$ws.Range("B22").Value = 0.6116819829310032
$ws.Range("C22").Value = 0.1204849657898848
$ws.Range("E22").Value = 0.7352026205710445
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.2636978349932377
$ws.Range("H22").Value = 0.443097488598184
$ws.Range("K22").Value = 0.6389998753467125
$ws.Range("O22").Value = 1.334255076996058

# Row 23
$ws.Range("B23").Value = 0.5902308406765542
$ws.Range("C23").Value = 0.117960752733282
$ws.Range("E23").Value = 0.7085544510077852
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.265078335137602
$ws.Range("H23").Value = 0.4452583060738746
$ws.Range("K23").Value = 0.6165100925632316
$ws.Range("O23").Value = 1.341615133149872

# Row 24
$ws.Range("B24").Value = 0.5088321942186269
$ws.Range("C24").Value = 0.1083785830861359
$ws.Range("E24").Value = 0.6078770647629739
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.270887048664207
$ws.Range("H24").Value = 0.4539263814929697
$ws.Range("K24").Value = 0.5310875659232295
$ws.Range("O24").Value = 1.371720845581194

# Row 25
$ws.Range("B25").Value = 0.420722901410528
$ws.Range("C25").Value = 0.09799931707446774
$ws.Range("E25").Value = 0.4997560556011678
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.2783799205204076
$ws.Range("H25").Value = 0.4643077370791389
$ws.Range("K25").Value = 0.4384469640354496
$ws.Range("O25").Value = 1.408926820221978

Write-Host "Updated 192 cells in Sheet1"
